# Add a "Git hub Link : <url>" line (with a live hyperlink on the URL)
# to the empty content placeholder on the "Github link" slide (slide 15).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(15)
$sh = $s.Shapes.Item("Content Placeholder 2")

$tr = $sh.TextFrame.TextRange

$label = "Git hub Link : "
$url   = "https://github.com/naveenip05/Predictive-Maintenance-of-Industrial-Machinery"

# Insert the full text into the (currently empty) paragraph. Using
# InsertAfter (rather than re-assigning .Text) keeps the paragraph's
# existing end-of-paragraph run properties intact.
$null = $tr.InsertAfter($label + $url)

# Label run.
$labelRange = $tr.Characters(1, $label.Length)
$labelRange.LanguageID = "en-IN"

# URL run, turned into a working external hyperlink.
$urlRange = $tr.Characters($label.Length + 1, $url.Length)
$urlRange.LanguageID = "en-IN"
$urlRange.ActionSettings.Item(1).Hyperlink.Address = $url
